$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.556.49'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.89%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.659.90'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.11%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.10%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '581.22'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.03%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.63'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.94%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.996'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.13%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.601'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.17%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '6.57'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.81%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.110'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.87%  '

# Row 11
$ws.Range('B11').Value = 'Cardano'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.382'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.00%  '

# Row 12
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.158'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.31%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.124.47'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.70%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.66'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +7.46%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '61.450.85'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.66%  '

# Row 16
$ws.Range('E16').Value = '  +3.59%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.673.96'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.31%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.64'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.38%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.79'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.68%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '355.84'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.43%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.94'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.36%  '

# Row 22
$ws.Range('E22').Value = '  +0.08%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.527'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.86%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.05'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.48%  '

# Row 25
$ws.Range('B25').Value = 'InternetComputer(DFINITY)'
$ws.Range('C25').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.61'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +7.11%  '

# Row 26
$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.164'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.53%  '

# Row 27
$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.995'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.26%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +7.32%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0827'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.11%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.85'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +7.63%  '

# Row 31
$ws.Range('B31').Value = 'USDe'
$ws.Range('C31').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.998'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.01%  '

# Row 32
$ws.Range('B32').Value = 'Monero'
$ws.Range('C32').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '166.98'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.84%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.13'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.99%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.73'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +11.11%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.11'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +13.61%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.35'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +8.93%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.73'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.40%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '346.73'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +11.90%  '

# Row 39
$ws.Range('E39').Value = '  +6.31%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.919'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +9.38%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '38.39'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.20%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.40'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +7.89%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0581'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.65%  '

# Row 44
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '21.26'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.29%  '

# Row 45
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '20.64'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.33%  '

# Row 46
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '135.37'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.30%  '

# Row 47
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.627'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.64%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0252'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.29%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.100'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.60%  '

# Row 50
$ws.Range('E50').Value = '  -0.26%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.103.77'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.35%  '
